# Generate Report for Handback
# Updates the handback-status workbook with refreshed timestamps and a
# status change (ht -> mt) on the "zh-cn" sheet, mirroring the rerun of
# the handback report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# first two files (both rows previously shared the same generate time).
$wsOverview.Range("G2").Value = "2016-08-21 12:15:44"
$wsOverview.Range("G3").Value = "2016-08-21 12:15:44"

# zh-cn sheet: Status column (E) flips from "ht" (human translation) to
# "mt" (machine translation) for the first two rows.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback
# DateTime (K) refreshed for the same two rows.
$wsZhCn.Range("H2").Value = "2016-08-21 12:15:40"
$wsZhCn.Range("H3").Value = "2016-08-21 12:15:40"
$wsZhCn.Range("K2").Value = "2016-08-21 12:15:57"
$wsZhCn.Range("K3").Value = "2016-08-21 12:15:57"

# de-de sheet: Correspond Handoff Datetime (H) mirrors the Overview
# sheet's "Latest HO Xliff Generate Date" value for the same two rows.
$wsDeDe.Range("H2").Value = "2016-08-21 12:15:44"
$wsDeDe.Range("H3").Value = "2016-08-21 12:15:44"

# de-de sheet: Correspond Handback DateTime (K) refreshed for the same
# two rows.
$wsDeDe.Range("K2").Value = "2016-08-21 12:16:07"
$wsDeDe.Range("K3").Value = "2016-08-21 12:16:07"
